$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F4 (485 -> 488) and F7 (659 -> 663)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 488
$ws1.Range("F7").Value = 663

# Sheet "全部类型" - same rows mirrored, update F4 (485 -> 488) and F7 (659 -> 663)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 488
$ws4.Range("F7").Value = 663
